$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-3 (Generation 0 and 1) to 7295
$ws.Range("C2:C3").Value = 7295

# Update rows 4-252 (Generation 2 through 250) to 7293
$ws.Range("C4:C252").Value = 7293
